# Update the division problems in the worksheet table.
# Values are stored as a 5-column table where the problem rows are
# row 1, 5, 9, 13, 17 (remaining rows are blank answer rows).

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# Row 1
$table.Cell(1,1).Range.Text = "29÷2="
$table.Cell(1,2).Range.Text = "93÷7="
$table.Cell(1,3).Range.Text = "14÷5="
$table.Cell(1,4).Range.Text = "30÷7="
$table.Cell(1,5).Range.Text = "33÷2="

# Row 5
$table.Cell(5,1).Range.Text = "23÷7="
$table.Cell(5,2).Range.Text = "44÷4="
$table.Cell(5,3).Range.Text = "31÷9="
$table.Cell(5,4).Range.Text = "27÷3="
$table.Cell(5,5).Range.Text = "10÷2="

# Row 9
$table.Cell(9,1).Range.Text = "20÷3="
$table.Cell(9,2).Range.Text = "78÷7="
$table.Cell(9,3).Range.Text = "40÷8="
$table.Cell(9,4).Range.Text = "87÷3="
$table.Cell(9,5).Range.Text = "84÷2="

# Row 13
$table.Cell(13,1).Range.Text = "95÷3="
$table.Cell(13,2).Range.Text = "28÷9="
$table.Cell(13,3).Range.Text = "43÷6="
$table.Cell(13,4).Range.Text = "56÷2="
$table.Cell(13,5).Range.Text = "83÷4="

# Row 17
$table.Cell(17,1).Range.Text = "25÷7="
$table.Cell(17,2).Range.Text = "54÷3="
$table.Cell(17,3).Range.Text = "60÷8="
$table.Cell(17,4).Range.Text = "42÷2="
$table.Cell(17,5).Range.Text = "21÷2="
